# إضافة عمود جديد 'Event ' إلى Card23
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Card23")

# Header cell M1, matching the style of the other header cells (bold/bordered/centered).
$ws.Cells.Item(1, 13).Value = "Event "
$ws.Cells.Item(1, 11).Copy()
$ws.Cells.Item(1, 13).PasteSpecial(-4122)

# Data cells M2:M12 - present but empty, same (default/plain) look as the rest of the table.
# A plain assignment of "" is treated as "clear the cell" (it would not persist at all), so
# nudge the cell into existing with a quote-prefixed empty string first, then paste the
# plain/default formatting from column A over it so the quote-prefix style doesn't show.
for ($r = 2; $r -le 12; $r++) {
    $ws.Cells.Item($r, 13).Value = "'"
    $ws.Cells.Item($r, 1).Copy()
    $ws.Cells.Item($r, 13).PasteSpecial(-4122)
}

$excel.CutCopyMode = 0
